# Applies the commit "Added some potential references to Discussion":
#  1. Extends the "Niche construction and cooperation" heading/line in the
#     Discussion outline with three citation markers.
#  2. Inserts three new Bibliography entries (Lehmann 2007, Platt and Bever
#     2009, Van Dyken and Wade 2012) in their correct alphabetical slots.

$d = $word.ActiveDocument

function Find-ParagraphIndex($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

function Insert-XmlIntoRange($range, $bodyXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1. "Niche construction and cooperation" -> add " *" and the three
#    citation runs.
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex("Niche construction and cooperation")
$p = $d.Paragraphs.Item($idx)
$runsXml = '<w:r><w:t xml:space="preserve">Niche construction and cooperation *</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">(Van Dyken and Wade, 2012)</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">*</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">(Lehmann, 2007)</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">*</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">(Platt and Bever, 2009)</w:t></w:r>'
$bodyXml = '<w:p>' + $runsXml + '</w:p>'
Insert-XmlIntoRange $p.Range $bodyXml

# ---------------------------------------------------------------------
# 2. New Bibliography entries.
# ---------------------------------------------------------------------

# -- Lehmann, 2007 (goes right after "Kuzdzal-Fick", before "McKinney") --
$idx = Find-ParagraphIndex("Kuzdzal-Fick")
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Item($idx + 1)
$bodyXml = '<w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Lehmann, L. 2007. The evolution of trans-generational altruism: Kin selection meets niche construction.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Journal of Evolutionary Biology</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">,</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">20</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">: 181&#8211;189. Blackwell Publishing Ltd.</w:t></w:r>' +
    '</w:p>'
Insert-XmlIntoRange $newp.Range $bodyXml

# -- Platt and Bever, 2009 (goes right after "Nowak", before "R Core Team") --
$idx = Find-ParagraphIndex("Nowak, M.A.")
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Item($idx + 1)
$bodyXml = '<w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Platt, T.G. and Bever, J.D. 2009. Kin competition and the evolution of cooperation.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Trends in Ecology &amp; Evolution</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">,</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">24</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">: 370&#8211;377.</w:t></w:r>' +
    '</w:p>'
Insert-XmlIntoRange $newp.Range $bodyXml

# -- Van Dyken and Wade, 2012 (goes right after "Us. 2015", before "Veelders") --
$idx = Find-ParagraphIndex("Us. 2015.")
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Item($idx + 1)
$bodyXml = '<w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Van Dyken, J.D. and Wade, M.J. 2012. Origins of altruism diversity II: Runaway coevolution of altruistic strategies via &#8220;reciprocal niche construction&#8221;.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Evolution</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">,</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">66</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">: 2498&#8211;2513.</w:t></w:r>' +
    '</w:p>'
Insert-XmlIntoRange $newp.Range $bodyXml

Write-Host "Edit complete."
